$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.107.10"
$ws.Range("E2").Value = "  -0.90%  "
$ws.Range("D3").Value = "1.665.00"
$ws.Range("E3").Value = "  -1.31%  "
$ws.Range("E4").Value = "  -0.84%  "
$ws.Range("D5").Value = "'209.47"
$ws.Range("E5").Value = "  -4.25%  "
$ws.Range("D6").Value = "'0.5174"
$ws.Range("E6").Value = "  -4.90%  "
$ws.Range("E7").Value = "  -0.78%  "
$ws.Range("D8").Value = "'0.2627"
$ws.Range("E8").Value = "  -4.43%  "
$ws.Range("D9").Value = "'0.06212"
$ws.Range("E9").Value = "  -3.62%  "
$ws.Range("E10").Value = "  -4.09%  "
$ws.Range("D11").Value = "'0.07483"
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("D12").Value = "1.711.17"
$ws.Range("E12").Value = "  +1.30%  "
$ws.Range("D13").Value = "'4.407"
$ws.Range("E13").Value = "  -2.59%  "
$ws.Range("D14").Value = "'0.5574"
$ws.Range("E14").Value = "  -4.12%  "
$ws.Range("D15").Value = "'65.88"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "'0.000007854"
$ws.Range("E16").Value = "  -6.29%  "
$ws.Range("D17").Value = "26.100.51"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "'4.767"
$ws.Range("E19").Value = "  -3.36%  "
$ws.Range("D20").Value = "'10.34"
$ws.Range("E20").Value = "  -5.57%  "
$ws.Range("D21").Value = "'185.88"
$ws.Range("E21").Value = "  -2.70%  "
$ws.Range("D22").Value = "'6.151"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").Value = "'147.59"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("D25").Value = "'0.1238"
$ws.Range("E25").Value = "  -6.35%  "
$ws.Range("D26").Value = "'7.533"
$ws.Range("E26").Value = "  -4.27%  "
$ws.Range("D27").Value = "'15.85"
$ws.Range("E27").Value = "  +0.74%  "
$ws.Range("D28").Value = "'0.06241"
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("E30").Value = "  -4.16%  "
$ws.Range("D31").Value = "'3.471"
$ws.Range("E31").Value = "  -2.79%  "
$ws.Range("D32").Value = "'3.408"
$ws.Range("E32").Value = "  -5.01%  "
$ws.Range("D33").Value = "'1.620"
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("D34").Value = "'0.9936"
$ws.Range("E34").Value = "  -4.53%  "
$ws.Range("E35").Value = "  -0.24%  "
$ws.Range("D36").Value = "'0.6008"
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("D37").Value = "'2.698"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "'6.121"
$ws.Range("E38").Value = "  -2.42%  "
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").Value = "1.072.73"
$ws.Range("E40").Value = "  -3.82%  "
$ws.Range("D41").Value = "'0.8601"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("D42").Value = "'1.004"
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("D43").Value = "'98.97"
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("D44").Value = "1.812.80"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("D45").Value = "'0.00000000110"
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("D46").Value = "'55.85"
$ws.Range("D47").Value = "'1.003"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").Value = "'0.05248"
$ws.Range("E48").Value = "  -0.43%  "
$ws.Range("D49").Value = "'7.897"
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").Value = "'5.900"
$ws.Range("E51").Value = "  -2.66%  "
